# Weekly data refresh: insert a new price record (2023-12-15, serial 45275)
# for "Repollo" / "Crespo record" / "Primera" at Macroferia Regional de Talca,
# pushing all existing rows from 491 downward by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 491; everything that was on/after row 491
# (through the former last row 593) shifts down to 492..594.
$ws.Rows.Item(491).Insert()

# Populate the newly inserted row with the new week's record.
$ws.Range("A491").Value = 5
$ws.Range("B491").Value = "Macroferia Regional de Talca"
$ws.Range("C491").Value = "Maule"
$ws.Range("D491").Value = 45275
$ws.Range("E491").Value = 7
$ws.Range("F491").Value = 100112006
$ws.Range("G491").Value = "Repollo"
$ws.Range("H491").Value = "Crespo record"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 3000
$ws.Range("K491").Value = 1200
$ws.Range("L491").Value = 1200
$ws.Range("M491").Value = 1200
$ws.Range("N491").Value = "$/unidad"
$ws.Range("O491").Value = "Región del Maule"
$ws.Range("P491").Value = 1200
$ws.Range("Q491").Value = 1
$ws.Range("R491").Value = "Hortaliza"
